$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 changes from the text "R40" to the text "1" (still stored as a
# shared string, not a number, and the cell keeps its original style).
# Writing a bare "1" through .Value lets Excel auto-coerce it to a real
# number, so instead compute it as text via TEXT(), then paste-special
# "values only" to freeze it as a literal string without disturbing the
# cell's existing number format/style.
$target = $ws.Range("B11")
$target.Formula = "=TEXT(1,""0"")"
$target.Copy() | Out-Null
$target.PasteSpecial(-4163)
$excel.CutCopyMode = $false
